$wb = $excel.ActiveWorkbook

# --- parameters sheet: move selection ---
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("D9").Select() | Out-Null

# --- coor sheet: selection becomes a range ---
$wsCoor = $wb.Worksheets.Item("coor")
$wsCoor.Range("D3:E9").Select() | Out-Null

# --- add new "test" sheet at the end of the tab strip ---
$wsTest = $wb.Worksheets.Add()
$wsTest.Name = "test"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTest.Move($null, $lastSheet)

# re-acquire the moved sheet object by name
$wsTest = $wb.Worksheets.Item("test")

# --- header row ---
$wsTest.Range("A1").Value = "hospital"
$wsTest.Range("B1").Value = "q1"
$wsTest.Range("C1").Value = "q2"

# --- population rows 1-7 (normal numeric style) ---
$popData = @(
    @(1, 5.25, 1.75),
    @(2, 22.5, 7.5),
    @(3, 12, 4),
    @(4, 6, 3),
    @(5, 15.75, 5.25),
    @(6, 11.25, 3.75),
    @(7, 14.25, 4.75)
)

$r = 2
foreach ($row in $popData) {
    $wsTest.Cells.Item($r, 1).Value = $row[0]
    $bCell = $wsTest.Cells.Item($r, 2)
    $bCell.Value = $row[1]
    $bCell.NumberFormat = "0.00"
    $cCell = $wsTest.Cells.Item($r, 3)
    $cCell.Value = $row[2]
    $cCell.NumberFormat = "0.00"
    $r = $r + 1
}

# --- extra highlighted rows 8-10 (size-50 population seed, repeated individuals) ---
$extraData = @(
    @(8, 5.25, 1.75),
    @(9, 22.5, 7.5),
    @(10, 12, 4)
)

foreach ($row in $extraData) {
    $wsTest.Cells.Item($r, 1).Value = $row[0]
    $bCell = $wsTest.Cells.Item($r, 2)
    $bCell.Value = $row[1]
    $bCell.NumberFormat = "0.00"
    $bCell.Interior.Color = 65535
    $cCell = $wsTest.Cells.Item($r, 3)
    $cCell.Value = $row[2]
    $cCell.NumberFormat = "0.00"
    $cCell.Interior.Color = 65535
    $r = $r + 1
}

# --- selection + activation on the new sheet ---
$wsTest.Range("F8:F9").Select() | Out-Null
$wsTest.Activate() | Out-Null
